$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New shared string text used for the Comprehension scores column (C10)
$newComprehension = "Gantz  (Text with visuals, Japanese, New):42; Harry Potter book 3 (Text-only, English, Familiar):33; ¿Qué sienten las plantas?  - CuriosaMente 260[https://www.youtube.com/watch?v=mGgnhpZ8d5g] (Audiovisual, Spanish, New):37; ¿Qué Pasaría si Minecraft Tuviera DEMASIADOS Bloques?[https://youtu.be/2pwjZfGOCTU] (Audiovisual, Spanish, New):39; "

# Row 9 gets an explicit custom row height (matches other data rows)
$ws.Rows.Item(9).RowHeight = 15.75

# New row 10 data
$ws.Range("A10").Value = 9

$ws.Range("B10").NumberFormat = "[h]:mm:ss"
$ws.Range("B10").Value = 2.2956597222222221

$ws.Range("C10").Value = $newComprehension

$ws.Range("D10").Value = $ws.Range("D9").Value()

# Update the active selection to C10, matching the saved workbook state
$ws.Range("C10").Select() | Out-Null
